# Update MLP model results (column G, "Valori stimati MLP") on sheets
# Caso1, Caso2 and Caso3 for rows 2-19.

$wb = $excel.ActiveWorkbook

$data = @{
    "Caso1" = @{
        2  = 1.01953387260437
        3  = 1.02214777469635
        4  = 1.019710302352905
        5  = 1.026696681976318
        6  = 1.029133081436157
        7  = 1.02318811416626
        8  = 1.022620916366577
        9  = 1.02508008480072
        10 = 1.023085594177246
        11 = 1.02351188659668
        12 = 1.022902965545654
        13 = 1.019252777099609
        14 = 1.023057222366333
        15 = 1.027686476707458
        16 = 1.020583033561707
        17 = 1.022072911262512
        18 = 1.022224307060242
        19 = 1.026463747024536
    }
    "Caso2" = @{
        2  = 1.021594643592834
        3  = 1.023833155632019
        4  = 1.021430492401123
        5  = 1.028966069221497
        6  = 1.030413389205933
        7  = 1.024636745452881
        8  = 1.024176955223083
        9  = 1.026232361793518
        10 = 1.024853467941284
        11 = 1.025332808494568
        12 = 1.024316906929016
        13 = 1.020912289619446
        14 = 1.02463972568512
        15 = 1.029813289642334
        16 = 1.02247416973114
        17 = 1.023364901542664
        18 = 1.023979306221008
        19 = 1.027941703796387
    }
    "Caso3" = @{
        2  = 1.020681381225586
        3  = 1.022453546524048
        4  = 1.020331621170044
        5  = 1.028133273124695
        6  = 1.029252529144287
        7  = 1.023064136505127
        8  = 1.023321270942688
        9  = 1.025023102760315
        10 = 1.02395486831665
        11 = 1.024209260940552
        12 = 1.022901773452759
        13 = 1.019897222518921
        14 = 1.02352774143219
        15 = 1.028404712677002
        16 = 1.021548628807068
        17 = 1.022204279899597
        18 = 1.0229572057724
        19 = 1.026838898658752
    }
}

foreach ($sheetName in $data.Keys) {
    $ws = $wb.Worksheets.Item($sheetName)
    $rows = $data[$sheetName]
    foreach ($r in $rows.Keys) {
        $ws.Cells.Item($r, 7).Value = $rows[$r]
    }
}
